$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record for Piña (Caramelo / Segunda) at "Macroferia Regional
# de Talca" needs to be inserted as row 206. Inserting the row shifts every
# existing record from row 206 down through row 260 down by one (to rows 207
# through 261), which matches the rest of the diff (each subsequent row's
# date/quality/price data becomes what used to be in the row above it).
$ws.Rows.Item(206).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A206").Value = 5
$ws.Range("B206").Value = "Macroferia Regional de Talca"
$ws.Range("C206").Value = "Maule"
$ws.Range("D206").Value = 44736
$ws.Range("E206").Value = 7
$ws.Range("F206").Value = "Fruta"
$ws.Range("G206").Value = 100108
$ws.Range("H206").Value = "Tropicales y subtropicales"
$ws.Range("I206").Value = 100108005
$ws.Range("J206").Value = "Piña"
$ws.Range("K206").Value = "Caramelo"
$ws.Range("L206").Value = "Segunda"
$ws.Range("M206").Value = 170
$ws.Range("N206").Value = 16000
$ws.Range("O206").Value = 16000
$ws.Range("P206").Value = 16000
$ws.Range("Q206").Value = "`$/caja 14 unidades"
$ws.Range("R206").Value = "Ecuador"
$ws.Range("S206").Value = 1143
$ws.Range("T206").Value = 14
